# Edit script: update loading_percent values for "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns (in order) that receive new values, mapped to their 1-based column index
$cols = @(2, 3, 4, 5, 7, 8, 9, 11, 13)   # B, C, D, E, G, H, I, K, M

# New values for rows 2..25 (one row per data record), in the same column order as $cols
$data = @(
            @(8.951501317152372, 5.655806928691322, 6.007734910445006, 12.08109677814984, 45.37897452412638, 18.00677494616066, 26.78271578696242, 9.705172977140139, 14.23580255800572),
            @(8.757089766904521, 5.445330853836549, 5.893435632429848, 11.84975282088854, 44.90887413912198, 17.97973445332703, 26.70904347305849, 9.578188036009605, 14.09607290731251),
            @(8.638696298478633, 5.313861239451921, 5.823978073490061, 11.70925077126039, 44.62802459879573, 17.96593432864122, 26.66805492204776, 9.502520859052773, 14.01377446825236),
            @(8.590778450961352, 5.259822361851697, 5.79589772242237, 11.65246548438505, 44.51564121010252, 17.96101725272683, 26.65242578788503, 9.472306319614269, 13.98115223965187),
            @(8.582844016182371, 5.250824275232329, 5.791249783403819, 11.643067144278, 44.49710774232643, 17.9602434952751, 26.64989561256453, 9.467327881285444, 13.97579162097346),
            @(8.638048621150778, 5.313134194190872, 5.823598408210271, 11.70848293167066, 44.6265004602924, 17.96586515234505, 26.66783978600303, 9.502110809081444, 14.01333076248602),
            @(8.88431501825122, 5.583755272436658, 5.968199930560149, 12.00105815189801, 45.21533075803879, 17.99687085566234, 26.75643621166295, 9.660936083150451, 14.18692029654482),
            @(9.371425323717201, 6.09282352385427, 6.255650660688172, 12.58346319626903, 46.42704286612565, 18.07981703430665, 26.96355600899423, 9.988777693808979, 14.5534008620206),
            @(9.72726639449659, 6.449092461515098, 6.466794655912919, 13.01193362962188, 47.34535986434565, 18.1541053652045, 27.13560373089815, 10.23706794799167, 14.83620912868332),
            @(9.887801318911475, 6.606521225100153, 6.562352737246082, 13.2060296154375, 47.76775453571002, 18.19075756650109, 27.21806800522, 10.35108700489234, 14.96728995429166),
            @(9.948329487197293, 6.665412595116027, 6.598428988882567, 13.27933605907153, 47.92824817833623, 18.20504338382778, 27.24988702529003, 10.39437192762447, 15.01723193132972),
            @(9.935306429493528, 6.652762282673215, 6.590664789059475, 13.26355798014022, 47.89366090711627, 18.20194867914688, 27.24300812270134, 10.38504564059352, 15.00646320582563),
            @(9.892786664265589, 6.611381050325359, 6.565323161902045, 13.21206488038428, 47.7809483059388, 18.19192475171446, 27.2206739919085, 10.35464616721347, 14.9713928655294),
            @(9.86670582777359, 6.585938083119236, 6.549785305189433, 13.18049648050553, 47.71197541325181, 18.18583759188401, 27.20707035677096, 10.33603838921538, 14.94994959111421),
            @(9.716741997845805, 6.438705780055473, 6.46053644886751, 12.99922597612177, 47.31783809684254, 18.15176718979938, 27.1302979661703, 10.22963416252874, 14.82768774989095),
            @(9.624349618069166, 6.347152802361404, 6.405630901061052, 12.88775785007759, 47.07714808132509, 18.1315950221857, 27.08426763119032, 10.16460013542816, 14.75327454707258),
            @(9.571085869026902, 6.294060117898427, 6.374006536030782, 12.82357172157822, 46.93915459466557, 18.12026184217115, 27.05818820916522, 10.12729773358609, 14.71070546615585),
            @(9.553032926250104, 6.276011170460895, 6.363292698470123, 12.80182930460763, 46.89251265934276, 18.11647100915884, 27.04942653727295, 10.11468695108463, 14.69633341950233),
            @(9.634198105123511, 6.356944103931323, 6.411480545954642, 12.89963189928113, 47.10272479968358, 18.13371454205935, 27.08912673897688, 10.17151274471167, 14.76117233270063),
            @(9.905283424270859, 6.623555769756378, 6.572769883240433, 13.22719549020165, 47.81404102685874, 18.19485803058342, 27.22721811363635, 10.36357264667902, 14.9816859686434),
            @(10.0808913939653, 6.793563682819468, 6.677526879881719, 13.44011686750933, 48.28202788467718, 18.23718534801998, 27.32091076454005, 10.48970618157513, 15.12755846638592),
            @(9.987331116688965, 6.703231835277422, 6.621688300160625, 13.32660701029592, 48.03201244563567, 18.21437955143573, 27.27059455459299, 10.42234521405436, 15.04955792078513),
            @(9.629746053889795, 6.352518877664795, 6.40883610057989, 12.89426395239892, 47.09116036433071, 18.13275548411651, 27.08692873906339, 10.16838728102551, 14.75760107984848),
            @(9.239687445789253, 5.957926378280864, 6.177725292961809, 12.42547312492298, 46.09381880310688, 18.05502100913925, 26.90399581477271, 9.898597352878818, 14.45170375847704)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $vals[$j]
    }
}
